# Update "想去人数" (interested-count) figures on the 展览 / 演出 / 全部类型
# sheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# 展览 (sheet1): rows 2,3,5,15,17,22,26,29,30,31,33,34
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14931
$ws1.Range("F3").Value = 18712
$ws1.Range("F5").Value = 123
$ws1.Range("F15").Value = 206
$ws1.Range("F17").Value = 1433
$ws1.Range("F22").Value = 7782
$ws1.Range("F26").Value = 1227
$ws1.Range("F29").Value = 110
$ws1.Range("F30").Value = 68
$ws1.Range("F31").Value = 164
$ws1.Range("F33").Value = 267
$ws1.Range("F34").Value = 5353

# 演出 (sheet2): row 3
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 4

# 全部类型 (sheet4): rows 2,3,5,15,17,23,27,29,32,33,34,36,37
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14931
$ws4.Range("F3").Value = 18712
$ws4.Range("F5").Value = 123
$ws4.Range("F15").Value = 206
$ws4.Range("F17").Value = 1433
$ws4.Range("F23").Value = 7782
$ws4.Range("F27").Value = 1227
$ws4.Range("F29").Value = 4
$ws4.Range("F32").Value = 110
$ws4.Range("F33").Value = 68
$ws4.Range("F34").Value = 164
$ws4.Range("F36").Value = 267
$ws4.Range("F37").Value = 5353
